$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 8995.64
$ws.Range("B9").Value = 8760
$ws.Range("C9").Value = 19.36
$ws.Range("D9").Value = 18.84
$ws.Range("E9").Value = $true
$ws.Range("F9").Value = -2.69
$ws.Range("G9").Value = 42612.672986111109
$ws.Range("G8").Copy()
$ws.Range("G9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H9").Value = $true
